$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save" — copy the header formatting from G1 (style index 1:
# bold, centered, bordered) onto H1, then set its text.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data column H2:H65 — the "Save" flag per row.
$saveValues = @(0,1,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,1,0,0,0,0,0,1,0,1,0,0,0,1,0,0,0,0,0,1,0,1,0,0,0,1,1,0,0,1,0,0,0,0,0,1,1,0,0,0,0,0,0,0,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
